$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the standalone "Meta description: ..." paragraph that originally
#    sat right under the "Play Fortune Rangers Free Slot - NetEnt Online
#    Game" H1 heading.
# ---------------------------------------------------------------------------
$metaIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        $metaIdx = $i
    }
}
if ($metaIdx -gt 0) {
    $d.Paragraphs.Item($metaIdx).Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph reading "Play Fortune Rangers Free Slot -
#    NetEnt Online Game" right before the final (italic) paragraph, which
#    immediately follows the "What we don't like" bullet list.
# ---------------------------------------------------------------------------
$anchorIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Might not be as appealing*") {
        $anchorIdx = $i
    }
}
$anchorPara = $d.Paragraphs.Item($anchorIdx)
$anchorRange = $anchorPara.Range
# A range collapsed to the interior of the anchor paragraph's text causes
# InsertXML to splice a whole new paragraph in before it, without disturbing
# the anchor paragraph's own runs (collapsing exactly on a paragraph
# boundary instead replaces that paragraph's content).
$insertPoint = $d.Range($anchorRange.Start + 1, $anchorRange.Start + 1)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fortune Rangers Free Slot - NetEnt Online Game</w:t></w:r></w:p>'
$insertPoint.InsertXML($newParaXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Swap the italic "Create a feature image for Fortune Rangers..." image
#    prompt text for the new meta-description copy, keeping the italic run
#    formatting intact.
# ---------------------------------------------------------------------------
$oldText = "Create a feature image for Fortune Rangers featuring a happy Maya warrior with glasses in cartoon style. The warrior should be holding a sword and wearing traditional Maya clothing with an Asian twist. The background should feature Chinese symbols and the Fortune Rangers logo. The overall style should be colorful and eye-catching to appeal to fans of Asian-inspired slot games and anime-style graphics."
$newText = "Read our Fortune Rangers slot review and play for free. Learn about Linked Reels feature, Asian theme, and more."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
